# Fixed issues with 81RF protective element
# Changed default xls parameters to disable 81x protections.
# Added goose messages for DER's cb's

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relays")
$ws.Activate()

# Disable the 81x (81RF) protections by changing the default relay
# parameters for every relay row (2-8):
#   81RFRP  [Hz/Sec] : column T / AH  -> 100
#   81RFDFP [Hz]     : column U / AI  -> 10
#   81RF Trip Delay [Sec] : column V / AJ -> 0.1
for ($r = 2; $r -le 8; $r++) {
    $ws.Range("T" + $r).Value = 100
    $ws.Range("U" + $r).Value = 10
    $ws.Range("V" + $r).Value = 0.1

    $ws.Range("AH" + $r).Value = 100
    $ws.Range("AI" + $r).Value = 10
    $ws.Range("AJ" + $r).Value = 0.1
}

# Reposition the window / selection the way it was left after editing the
# new AH:AJ parameter block.
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AH2:AJ8").Select()
